$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")
$ws.Columns("AS:AS").Delete()
